# Doing Updates for Financials
# Insert a new first data column (fiscal year 2018 / 31-Dec-2018) into the
# "IAC" income statement / balance sheet / cash flow statement tables.
# This shifts the existing data columns D:K one position to the right
# (becoming E:L) and populates the new column D with the FY2018 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a new column before column D - this shifts D:K -> E:L and
#    keeps all existing values/formatting intact on the shifted columns.
$ws.Columns("D").Insert()

# 2. The freshly inserted column D has no real formatting yet (it just
#    inherited the blank default). Copy number formatting/style from the
#    column immediately to its right (E, which used to be D) for every
#    row that actually participates in one of the three data tables.
#    (Rows that are pure section headers - 5, 6, 37, 79 - never had a D
#    cell and must stay that way.)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# 3. Populate the new column D with the FY2018 (31-Dec-2018) values.

# -- Income Statement --------------------------------------------------
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 4262900
$ws.Range("D9").Value2 = 911100
$ws.Range("D10").Value2 = 3351700
$ws.Range("D12").Value2 = 309300
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = -91800
$ws.Range("D15").Value2 = 155000
$ws.Range("D17").Value2 = 3577200
$ws.Range("D18").Value2 = 685700
$ws.Range("D20").Value2 = 185100
$ws.Range("D21").Value2 = 1054600
$ws.Range("D22").Value2 = 109300
$ws.Range("D23").Value2 = 761600
$ws.Range("D24").Value2 = 13000
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 748600
$ws.Range("D27").Value2 = 617800
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 9200
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -185100
$ws.Range("D33").Value2 = 627000
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 627000

# -- Balance Sheet ------------------------------------------------------
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 2131600
$ws.Range("D42").Value2 = 123700
$ws.Range("D43").Value2 = 289300
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 218100
$ws.Range("D46").Value2 = 2762700
$ws.Range("D47").Value2 = 235100
$ws.Range("D48").Value2 = 318800
$ws.Range("D49").Value2 = 3358300
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 199700
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 6874600
$ws.Range("D57").Value2 = 74900
$ws.Range("D58").Value2 = 13800
$ws.Range("D59").Value2 = 794900
$ws.Range("D60").Value2 = 883600
$ws.Range("D61").Value2 = 2245500
$ws.Range("D62").Value2 = 128000
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 4031500
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 1258800
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 2843100
$ws.Range("D77").Value2 = 0

# -- Cash Flow Statement --------------------------------------------------
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 627000
$ws.Range("D83").Value2 = 183800
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 988100
$ws.Range("D91").Value2 = -85600
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -173400
$ws.Range("D96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -312800
$ws.Range("D101").Value2 = -1900
$ws.Range("D102").Value2 = 500000

# 4. A couple of cells in the shifted-over columns were not merely shifted
#    values from the old layout - the author touched them at the same
#    time, so they need an explicit correction after the shift:
$ws.Range("F89").Value2 = 344200
$ws.Range("E94").Value2 = 42000
